# "multiply matrices transposed (results)"
#
# The `direct` vs `transposed` matrix-multiply benchmark grew two more
# matrix sizes (1400 and 1500) and every existing measurement was
# re-run (the whole "direct"/"transposed" column shifted - hence every
# value changing even for the sizes that already existed), so rows
# 24-30 get new values and rows 31-32 are brand new.  The "true
# sharing" / "false sharing" chart above it is untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- benchmark data: matrix size (A) / direct (B) / transposed (C) ----
$data = @(
    @(24, 700,  669,  344),
    @(25, 800,  807,  478),
    @(26, 900,  1357, 700),
    @(27, 1000, 1788, 1033),
    @(28, 1100, 3104, 1348),
    @(29, 1200, 4117, 1755),
    @(30, 1300, 6722, 2279),
    @(31, 1400, 9235, 2838),
    @(32, 1500, 14018, 3431)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
}

# --- grow the "direct vs transposed" chart's series to the new range --
$chart2 = $ws.ChartObjects(2).Chart
$directSeries = $chart2.SeriesCollection(1)
$directSeries.Formula = "=SERIES(benchmarks!`$B`$23,,benchmarks!`$B`$24:`$B`$32,1)"
$transposedSeries = $chart2.SeriesCollection(2)
$transposedSeries.Formula = "=SERIES(benchmarks!`$C`$23,,benchmarks!`$C`$24:`$C`$32,2)"

# --- re-flow the two chart frames on the sheet (moved/resized) --------
$chart1Obj = $ws.ChartObjects(1)
$chart1Obj.Left = 379.984251968504
$chart1Obj.Top = 9.04251968503937
$chart1Obj.Width = 447.137007874016
$chart1Obj.Height = 251.489763779528

$chart2Obj = $ws.ChartObjects(2)
$chart2Obj.Left = 319.492913385827
$chart2Obj.Top = 289.700787401575
$chart2Obj.Width = 453.514960629921
$chart2Obj.Height = 255.089763779528

Write-Host "Applied matrix transpose benchmark update"
